# Updated cryptos list on Tue Nov  7 02:46:39 UTC 2023 with GitHub Actions
# Refreshes price/volume columns (and two re-ranked coin rows) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.184.00"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").Value = "1.900.37"
$ws.Range("E3").Value = "  -0.10%  "

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.692"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.20"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.362"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0758"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0984"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.10%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.59%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.177.22"
$ws.Range("E14").Value = "  -0.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.737"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.66%  "

$ws.Range("D17").Value = "1.876.28"
$ws.Range("E17").Value = "  -1.45%  "

$ws.Range("D18").Value = "35.173.41"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.03%  "

$ws.Range("D20").Value = "0.0₃0836"
$ws.Range("E20").Value = "  +1.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "242.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.66%  "

$ws.Range("E23").Value = "  +4.19%  "

$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("E25").Value = "  +4.54%  "

$ws.Range("E26").Value = "  -0.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.63%  "

$ws.Range("E30").Value = "  -1.32%  "

$ws.Range("D31").Value = "4.127.90"
$ws.Range("E31").Value = "  -0.56%  "

$ws.Range("E32").Value = "  +18.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0606"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.18%  "

$ws.Range("E34").Value = "  +3.32%  "

$ws.Range("E35").Value = "  +16.73%  "

$ws.Range("E36").Value = "  +1.51%  "

$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.857"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -13.21%  "

$ws.Range("E39").Value = "  -2.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0214"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.03%  "

$ws.Range("E43").Value = "  +0.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0648"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.58%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.330.43"
$ws.Range("E46").Value = "  -0.89%  "

$ws.Range("E47").Value = "  +0.65%  "

$ws.Range("E48").Value = "  -1.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.86%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.55%  "

